# Update to respiration routine: add the new RespParam row used by the
# updated Resp routine (general housekeeping / parameter table update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New parameter row appended right after the existing Min/Resp parameters.
$ws.Range("A22").Value = "RespParam"
$ws.Range("B22").Value = 0.001
$ws.Range("C22").Value = "unitless"

# Leave the selection where the author left it after entering the new value.
$ws.Range("D21").Select()
